$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "BL"
$ws.Range("G1").Value = "Operating Freq"

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
}

$ws.Columns.Item(7).AutoFit() | Out-Null

$ws.Range("I9").Select() | Out-Null
